# Daily attendance processing - 2026-01-06 10:05:09
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# lists the users who recorded each attendance session. For a specific
# set of rows, the order of the two names in that comma-separated list
# is flipped from "dnasr281@gmail.com, System" to
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Rows (in column G) whose "Recorded By" value needs to be reordered.
$rows = @(2, 3, 4, 5, 6, 7, 16, 17, 22, 23, 37, 38, 43, 44, 58, 59, 64, 65, 79, 80, 85, 86, 87, 88, 89, 90, 99, 100, 105, 106, 107, 108, 109, 110, 119, 120, 125, 126, 127, 128, 129, 130, 139, 140, 145, 146, 147, 148, 149, 150, 159, 160, 165, 166, 167, 168, 169, 170, 179, 180, 185, 186, 200, 201, 206, 207, 221, 222, 227, 228, 242, 243)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
